$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new row as text-producing formulas first so Excel does not
# auto-convert the look-alike numeric/date strings into real numbers/dates,
# then convert them to static values (Copy / Paste Special - Values) so the
# cells end up holding plain text, matching the source data export.
$ws.Range("A2").Formula = '="07/07/2023"'
$ws.Range("B2").Formula = '="1000.00"'
$ws.Range("C2").Formula = '="1000.00"'
$ws.Range("D2").Formula = '="1000.00"'
$ws.Range("E2").Formula = '="1000.00"'
$ws.Range("F2").Formula = '="0.00"'
$ws.Range("G2").Formula = '="100.00"'

$rng = $ws.Range("A2:G2")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0
